# Inserts one new daily price record for "Pepino ensalada" at row 41
# (Vega Monumental Concepción), pushing the existing rows 41-100 down to
# 42-101 — matches the weekly-refresh pattern described in the commit
# message ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 41..100 down to 42..101, leaving a blank row 41 to fill in.
$ws.Rows.Item(41).Insert()

$ws.Cells.Item(41, 1).Value  = 11
$ws.Cells.Item(41, 2).Value  = 'Vega Monumental Concepción'
$ws.Cells.Item(41, 3).Value  = 'Bíobío'
$ws.Cells.Item(41, 4).Value  = 44580
$ws.Cells.Item(41, 5).Value  = 8
$ws.Cells.Item(41, 6).Value  = 100112043
$ws.Cells.Item(41, 7).Value  = 'Pepino ensalada'
$ws.Cells.Item(41, 8).Value  = 'Sin especificar'
$ws.Cells.Item(41, 9).Value  = 'Primera'
$ws.Cells.Item(41, 10).Value = 150
$ws.Cells.Item(41, 11).Value = 11000
$ws.Cells.Item(41, 12).Value = 12000
$ws.Cells.Item(41, 13).Value = 11667
$ws.Cells.Item(41, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(41, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(41, 16).Value = 194
$ws.Cells.Item(41, 17).Value = 60
$ws.Cells.Item(41, 18).Value = 'Hortaliza'
